$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix typo on existing row 15: "Realisation and Desing" -> "Realisation and Design" ---
$ws.Range("E15").Value = "Realisation and Design"

# --- Copy the formatting of row 15 down into the 4 new log rows (16-19) ---
$ws.Range("A15:F15").Copy()
$ws.Range("A16:F19").PasteSpecial(-4122)   # xlPasteFormats

# Row 16: 10-3-2010
$ws.Range("A16").Value = 40247
$ws.Range("B16").Value = 0.39583333333333331
$ws.Range("C16").Value = 0.64583333333333337
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = "Realisation and Design"
$ws.Range("F16").Value = "Design multiple kinds of surfaces, enemy implementation"

# Row 17: 11-3-2010
$ws.Range("A17").Value = 40248
$ws.Range("B17").Value = 0.4375
$ws.Range("C17").Value = 0.64583333333333337
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = "Realisation and Design"
$ws.Range("F17").Value = "Enemy factory"

# Row 18: 15-3-2010 (hours logged as "4.5" text, like the earlier 3.5/6.5 entries)
$ws.Range("A18").Value = 40252
$ws.Range("B18").Value = 0.45833333333333331
$ws.Range("C18").Value = 0.60416666666666663
$ws.Range("Z1").Formula = '="4.5"'
$ws.Range("Z1").Copy()
$ws.Range("D18").PasteSpecial(-4163)       # xlPasteValues - keeps text type, keeps D18's existing style
$ws.Range("Z1").Clear()
$ws.Range("E18").Value = "Design"
$ws.Range("F18").Value = "Tiles for the snow level"

# Row 19: 16-3-2010 (hours logged as "4.75" text)
$ws.Range("A19").Value = 40253
$ws.Range("B19").Value = 0.42708333333333331
$ws.Range("C19").Value = 0.625
$ws.Range("Z1").Formula = '="4.75"'
$ws.Range("Z1").Copy()
$ws.Range("D19").PasteSpecial(-4163)       # xlPasteValues
$ws.Range("Z1").Clear()
$ws.Range("E19").Value = "Realisation"
$ws.Range("F19").Value = "Snowlevel implementation, all tiles and some surfaces"

# --- Pre-format a big block of empty rows below (D20:D70), matching D-column's right-aligned style ---
$ws.Range("D20:D70").HorizontalAlignment = -4152   # xlRight

# --- Update the view state: scrolled down a bit, selection on F21 ---
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("F21").Select()
